$wb = $excel.ActiveWorkbook

# --- Update the status text everywhere it appears ("Ready for handoff" -> "Handed back: in sync with en-US") ---
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# --- zh-cn sheet: fill in Latest Target File (I) / Latest Handback File (J) ---
$zhRow2TargetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/900f8d72356ee96e9c3ad74adf206353ceb93fdf/e2e/0285d837-b6da-4c87-86d9-8789404b5f20.md"
$zhRow3TargetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/900f8d72356ee96e9c3ad74adf206353ceb93fdf/e2e/55ed2204-3c8a-4d10-a59c-ae1e9528a8ea.md"

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $zhRow2TargetUrl, "", "", "0285d837-b6da-4c87-86d9-8789404b5f20.md")
$wsZh.Range("J2").Value = $wsZh.Range("G2").Value2

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $zhRow3TargetUrl, "", "", "55ed2204-3c8a-4d10-a59c-ae1e9528a8ea.md")
$wsZh.Range("J3").Value = $wsZh.Range("G3").Value2

# --- de-de sheet: fill in Latest Target File (I) / Latest Handback File (J) / Latest Handback DateTime (K) ---
$deRow2TargetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/900f8d72356ee96e9c3ad74adf206353ceb93fdf/e2e/0285d837-b6da-4c87-86d9-8789404b5f20.md"
$deRow3TargetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/900f8d72356ee96e9c3ad74adf206353ceb93fdf/e2e/55ed2204-3c8a-4d10-a59c-ae1e9528a8ea.md"

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $deRow2TargetUrl, "", "", "0285d837-b6da-4c87-86d9-8789404b5f20.md")
$wsDe.Range("J2").Value = $wsDe.Range("G2").Value2
$wsDe.Range("K2").Value = "2016-08-24 10:29:48"

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $deRow3TargetUrl, "", "", "55ed2204-3c8a-4d10-a59c-ae1e9528a8ea.md")
$wsDe.Range("J3").Value = $wsDe.Range("G3").Value2
$wsDe.Range("K3").Value = "2016-08-24 10:29:48"

Write-Output "done"
